$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organizations")
Write-Host $ws.Name
